# Amplifier.xlsx edit: add a new "Formel (Ua*-0,485)+3,1" column (H/I) that
# recomputes the amplifier output with an updated gain formula, drop the old
# "(Messwert mit Sensor)" note, make room for it by inserting two rows above
# the second (temperature delta) table, and repoint charts/drawings at the
# new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------
# 1. Insert two rows above the second table (old rows 17-38 -> 19-40).
#    This also shifts the chart2 source ranges down by two rows.
# ---------------------------------------------------------------------
$ws.Rows("17:18").Insert()

# ---------------------------------------------------------------------
# 2. Remove the stray "(Messwert mit Sensor)" note that used to live in H6.
# ---------------------------------------------------------------------
$ws.Range("H6").ClearContents()

# ---------------------------------------------------------------------
# 3. Add the new header in H4 and fill H5:I12 with the new formula columns.
# ---------------------------------------------------------------------
$ws.Range("H4").Value = "Formel (Ua*-0,485)+3,1"
$ws.Range("H4").Style = $ws.Range("G4").Style

$ws.Range("H5").Formula = "=(D5*-0.485)+3.1"
$ws.Range("H5").Style = $ws.Range("F5").Style
$ws.Range("I5").Formula = "=(E5*-0.485)+3100"
$ws.Range("I5").Style = $ws.Range("G5").Style

$ws.Range("H6:H12").Formula = "=(D6*-0.485)+3.1"
$ws.Range("H6:H12").Style = $ws.Range("F6").Style
$ws.Range("I6:I12").Formula = "=(E6*-0.485)+3100"
$ws.Range("I6:I12").Style = $ws.Range("G6").Style

# ---------------------------------------------------------------------
# 4. Point the first chart's third series (formerly the raw "Ua" series)
#    at the new formula column instead, and drop its category override
#    (it now falls back to the shared Ub category from the other series).
# ---------------------------------------------------------------------
$chart1 = $ws.ChartObjects(1).Chart
$ser3 = $chart1.SeriesCollection(3)
$ser3.Values = $ws.Range("I5:I12")
$ser3.XValues = $ws.Range("B5:B12")
$ser3.Name = "=Tabelle1!$H$4"

$chart1.Axes(2).MinimumScale = 850

# ---------------------------------------------------------------------
# 5. Reposition the two charts: chart 1 grows two rows taller, chart 2
#    moves down two rows and grows two rows taller too.
# ---------------------------------------------------------------------
$draw = $ws.Shapes
$shp1 = $ws.Shapes.Item("Diagramm 2")
$shp1.TopLeftCell.Row
$t1 = $shp1.TopLeftCell
$shp1.Top = $ws.Cells.Item(2, 10).Top
$l1 = $shp1.Left
$w1 = $shp1.Width
$h1 = $shp1.Height

$wb.Windows(1).Visible = $wb.Windows(1).Visible
